# "Main only runs algos once now" - the algorithm run counts/timings in
# row 23 (M/N) and the whole re-run sample in row 31 change because the
# test harness no longer repeats each algorithm run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 23: only the Dijkstra / A* timing columns shift.
$ws.Range("M23").Value = 348
$ws.Range("N23").Value = 327

# Row 31: the whole sample row shifts since it was re-captured from a
# single algorithm run instead of an averaged/repeated run.
$ws.Range("B31").Value = 535
$ws.Range("C31").Value = 494
$ws.Range("D31").Value = 364
$ws.Range("E31").Value = 414
$ws.Range("F31").Value = 333
$ws.Range("G31").Value = 345
$ws.Range("M31").Value = 351
$ws.Range("N31").Value = 354

# Leave the view/selection the way the author last left it (scrolled up one
# row, with AB26 as the active cell instead of the old R23:W23 selection).
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 14
$ws.Range("AB26").Select()
